# Slide 4 ("Code 1/3"): re-crop the code-screenshot picture.
#
# The picture (p:pic "Picture 3", r:embed rId2, native size 518x599px ==
# 388.5pt x 449.25pt) gets its left edge cropped by 4.167% (srcRect l="4167")
# while its right edge on the slide stays put - i.e. the shape's Left grows
# and its Width shrinks by exactly the cropped-away slice
# (388.5pt * 0.04167 = 16.188795pt = 205,632/12700 EMU -> 216024 EMU),
# so the visible picture shifts left by that fraction without otherwise
# moving or resizing.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item("Picture 3")

# Crop 4.167% off the left of the source image.
$shp.PictureFormat.CropLeft = 16.188795

# Keep the right edge (Left+Width) fixed at 638.1763in.. err, pt, by moving
# Left right and shrinking Width by the same amount that was cropped away.
$shp.Left = 229.59181102362206
$shp.Width = 391.22456692913386
